$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Start period", "End period" and "Apartment" columns (D:F).
# This shifts "Payment date" (was G) to D and "Note" (was H) to E.
$ws.Range("D1:F1").EntireColumn.Delete()

# The "Revenue type" header becomes "Invoice".
$ws.Range("B1").Value = "Invoice"

# Selection moves to A2.
$ws.Range("A2").Select()
